$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to Text format, write the literal string, then strip
    # the temporary format again so the final style matches an untouched cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '23.552.06'
$ws.Range("E2").Value = '  +1.56%  '
Set-TextValue $ws.Range("D3") '1.656.29'
$ws.Range("E3").Value = '  +2.94%  '
Set-TextValue $ws.Range("D4") '1.000'
$ws.Range("E4").Value = '  -0.52%  '
Set-TextValue $ws.Range("D5") '0.9997'
$ws.Range("E5").Value = '  -0.43%  '
Set-TextValue $ws.Range("D6") '302.37'
$ws.Range("E6").Value = '  +0.08%  '
Set-TextValue $ws.Range("D7") '0.3835'
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range("D8") '0.3599'
$ws.Range("E8").Value = '  +2.54%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D9") '51.08'
$ws.Range("E9").Value = '  -1.70%  '
Set-TextValue $ws.Range("D10") '0.08207'
$ws.Range("E10").Value = '  +1.95%  '
Set-TextValue $ws.Range("D11") '1.237'
$ws.Range("E11").Value = '  +3.89%  '
Set-TextValue $ws.Range("D12") '1.000'
$ws.Range("E12").Value = '  -0.54%  '
Set-TextValue $ws.Range("D13") '22.42'
$ws.Range("E13").Value = '  +2.25%  '
Set-TextValue $ws.Range("D14") '6.481'
$ws.Range("E14").Value = '  +2.54%  '
Set-TextValue $ws.Range("D15") '7.520'
$ws.Range("E15").Value = '  +4.10%  '
Set-TextValue $ws.Range("D16") '0.00001225'
$ws.Range("E16").Value = '  +1.45%  '
Set-TextValue $ws.Range("D17") '1.652.56'
$ws.Range("E17").Value = '  +4.05%  '
$ws.Range("E18").Value = '  +3.67%  '
Set-TextValue $ws.Range("D19") '0.06984'
$ws.Range("E19").Value = '  +1.11%  '
Set-TextValue $ws.Range("D20") '6.800'
$ws.Range("E20").Value = '  +5.40%  '
$ws.Range("E21").Value = '  +2.99%  '
Set-TextValue $ws.Range("D22") '0.9994'
$ws.Range("E22").Value = '  -0.44%  '
$ws.Range("E23").Value = '  +4.02%  '
Set-TextValue $ws.Range("D24") '23.568.20'
$ws.Range("E24").Value = '  +1.69%  '
Set-TextValue $ws.Range("D25") '2.518'
$ws.Range("E25").Value = '  -1.05%  '
Set-TextValue $ws.Range("D26") '3.024'
$ws.Range("E26").Value = '  -0.67%  '
$ws.Range("E27").Value = '  +2.49%  '
Set-TextValue $ws.Range("D28") '152.62'
$ws.Range("E28").Value = '  +1.14%  '
Set-TextValue $ws.Range("D29") '5.244'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("E30").Value = '  +2.09%  '
Set-TextValue $ws.Range("D31") '1.832.84'
$ws.Range("E31").Value = '  +3.63%  '
Set-TextValue $ws.Range("D32") '7.162'
$ws.Range("E32").Value = '  +12.18%  '
Set-TextValue $ws.Range("D33") '2.247'
$ws.Range("E33").Value = '  +6.10%  '
Set-TextValue $ws.Range("D34") '12.03'
$ws.Range("E34").Value = '  +6.03%  '
Set-TextValue $ws.Range("D35") '1.062'
$ws.Range("E35").Value = '  +0.54%  '
Set-TextValue $ws.Range("D36") '0.02801'
$ws.Range("E36").Value = '  +4.03%  '
Set-TextValue $ws.Range("D37") '6.117'
$ws.Range("E37").Value = '  +5.51%  '
$ws.Range("E38").Value = '  +2.20%  '
Set-TextValue $ws.Range("D39") '0.08778'
$ws.Range("E39").Value = '  +1.39%  '
Set-TextValue $ws.Range("D40") '0.07001'
$ws.Range("E40").Value = '  +1.99%  '
$ws.Range("E41").Value = '  +10.87%  '
Set-TextValue $ws.Range("D42") '0.6997'
$ws.Range("E42").Value = '  +2.82%  '
Set-TextValue $ws.Range("D43") '1.333'
$ws.Range("E43").Value = '  +2.12%  '
Set-TextValue $ws.Range("D44") '15.94'
$ws.Range("E44").Value = '  +5.03%  '
Set-TextValue $ws.Range("D45") '0.6534'
$ws.Range("E45").Value = '  +4.69%  '
Set-TextValue $ws.Range("D46") '1.000'
$ws.Range("E46").Value = '  -0.36%  '
Set-TextValue $ws.Range("D47") '2.305'
$ws.Range("E47").Value = '  +3.16%  '
$ws.Range("E48").Value = '  +0.54%  '
Set-TextValue $ws.Range("D49") '0.07902'
$ws.Range("E49").Value = '  +0.73%  '
Set-TextValue $ws.Range("D50") '128.39'
$ws.Range("E50").Value = '  +0.64%  '
Set-TextValue $ws.Range("D51") '1.190'
$ws.Range("E51").Value = '  +2.68%  '
